# Add team record columns (Wins, Losses, Ties) to the right of the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used data row (row 1 is the header, data runs 2..52).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New header cells in row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) used by
# the other header cells such as AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill Wins/Losses/Ties for every data row with the team's season record.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 51   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 111  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
